# "update to the ppt done"
# Update the bullet about the touch screen setup on slide 4 of the
# project-status content placeholder to also mention that the ROS
# installation is done.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)            # "Content Placeholder 8"
$tr = $shp.TextFrame.TextRange

$oldText = "Touch screen setup on raspberry pi."
$newText = "Touch screen setup on raspberry pi and ROS installation done."

# Locate the paragraph that currently holds the bullet we need to edit.
# Note: TextRange.Text for a single paragraph includes a trailing "`r"
# paragraph-end marker, so trim it before comparing.
$paragraphs = $tr.Paragraphs()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.TrimEnd("`r") -eq $oldText) {
        # Replace the text through a Characters() range spanning the whole
        # run so PowerPoint updates the existing run's text in place
        # (keeping a single run with its original rPr) instead of
        # diff-splitting it into multiple runs.
        $chars = $para.Characters(1, $para.Length)
        $chars.Text = $newText
        break
    }
}
